$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Move the "栈指针ESP" label (and its vertically-centered style) from A10 up to A9
$ws.Range("A9").Value = $ws.Range("A10").Value()
$ws.Range("A9").VerticalAlignment = -4108

# Fully clear old A10 cell (value + formatting) so it no longer exists in the sheet
$ws.Range("A10").Clear()

# Update the active selection to A9
$ws.Range("A9").Select()
